$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before FB (01-oct.) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Insert a new entire column at FB; everything from FB.. shifts right to FC..GG
$ws1.Range("FB1").EntireColumn.Insert()

# New column header continues the December date sequence (FA1 = "30-dec")
$ws1.Range("FB1").Value2 = "31-dec"

# New column has no data for this date yet -> "-" placeholder, like the other
# not-yet-happened days in the sheet
$ws1.Range("FB2:FB25").Value2 = "-"

# --- Sheet "Gaz": append the new daily price row ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A186").NumberFormat = "@"
$ws2.Range("A186").Value2 = "2025-12-29"
$ws2.Range("A186").ClearFormats()
$ws2.Range("B186").Value2 = 28.105

# --- Sheet "CO2": append the new daily price row ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A183").NumberFormat = "@"
$ws3.Range("A183").Value2 = "2025-12-29"
$ws3.Range("A183").ClearFormats()
$ws3.Range("B183").Value2 = 85.73999999999999
